$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell H1: "Save", styled like the other header cells (B1:G1)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# Data values for H2:H8
$values = @(0, 1, 0, 1, 0, 0, 0)
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 8).Value = $values[$i]
}
